# Avances.xlsx data refresh - "Se estabiliza la transacciòn de avances solicitada"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# --- Preserve the quote-prefixed text style used by the numeroTarjeta (M) column ---
# M2:M8 use a cell style with a leading quotePrefix; setting .Value directly on those
# cells would silently switch them to a plain text style, so we stash that exact
# formatting in an unused scratch cell first and paste it back after the values change.
$ws.Range("M2").Copy() | Out-Null
$ws.Range("Z1").PasteSpecial(-4122) | Out-Null

# --- usuario (D) column: testing10 -> avancestdc8 ---
$ws.Range("D2:D6").Value = "avancestdc8"

# --- numeroTarjeta (M) column: new masked card numbers ---
$ws.Range("M2").Value = "*2442"
$ws.Range("M3").Value = "*6716"
$ws.Range("M4").Value = "*2442"
$ws.Range("M5").Value = "*6716"
$ws.Range("M6").Value = "*8078"

# Restore the quotePrefix formatting on the M column after rewriting the values.
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("M2:M6").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Clear() | Out-Null

# --- codigoSeguridad (N) column: new security codes ---
$ws.Range("N2").Value = 314
$ws.Range("N3").Value = 9639
$ws.Range("N4").Value = 314
$ws.Range("N5").Value = 9639
$ws.Range("N6").Value = 9639

# --- montoAvance (O) column: updated amounts ---
$ws.Range("O2").Value = 250000
$ws.Range("O3").Value = 150000
$ws.Range("O4").Value = 20000
$ws.Range("O5").Value = 1300000
$ws.Range("O6").Value = 1187500

# --- numeroTarjetaDestino (P) column: new destination card ---
$ws.Range("P2:P6").Value = "406-110080-05"

# --- Remove the two trailing sample rows (old rows 7 and 8) ---
$ws.Rows("7:8").Delete()

# --- Refresh the cross-sheet list data validation so it only spans the remaining rows ---
$ws.Range("G2:J6").Validation.Delete()
$ws.Range("G2:J6").Validation.Add(3, 1, 1, "=Listas!`$A`$2:`$A`$3")

# --- Update the visible selection/scroll position to match the edited view ---
$ws.Activate() | Out-Null
$ws.Range("I1").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 9
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("O7").Select() | Out-Null
